$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F393").Value = 308148  # was 308126

$ws.Range("F394").Value = 166181  # was 166172

$ws.Range("F395").Value = 750844  # was 750792

$ws.Range("F398").Value = 298498  # was 298177

$ws.Range("F399").Value = 200365  # was 200436
$ws.Range("G399").Value = 965  # was 966

$ws.Range("F400").Value = 150767  # was 150762

$ws.Range("F401").Value = 273360  # was 273355

$ws.Range("F402").Value = 717945  # was 716948
$ws.Range("G402").Value = 1387  # was 1388

$ws.Range("F403").Value = 351936  # was 351927

$ws.Range("F405").Value = 173936  # was 173903
$ws.Range("G405").Value = 694  # was 693

$ws.Range("F407").Value = 158126  # was 158120

$ws.Range("F408").Value = 303658  # was 303656

$ws.Range("F409").Value = 704350  # was 703004

$ws.Range("F410").Value = 363442  # was 363434

$ws.Range("F411").Value = 224997  # was 225001

$ws.Range("F412").Value = 175809  # was 175804

$ws.Range("F413").Value = 148664  # was 148919
$ws.Range("G413").Value = 658  # was 659

$ws.Range("F414").Value = 146998  # was 146763
$ws.Range("G414").Value = 558  # was 557

$ws.Range("F415").Value = 304014  # was 304867
$ws.Range("G415").Value = 692  # was 693

$ws.Range("F416").Value = 660106  # was 658759
$ws.Range("G416").Value = 926  # was 923

$ws.Range("F417").Value = 332536  # was 332534

$ws.Range("F418").Value = 200487  # was 200363

$ws.Range("F419").Value = 147466  # was 147463

$ws.Range("F420").Value = 136772  # was 136598

$ws.Range("F422").Value = 294118  # was 293474

$ws.Range("F424").Value = 255482  # was 255004

$ws.Range("F426").Value = 105019  # was 104543

$ws.Range("F427").Value = 89163  # was 89031
$ws.Range("G427").Value = 359  # was 358

$ws.Range("F428").Value = 100095  # was 99209
$ws.Range("G428").Value = 382  # was 374

$ws.Range("F429").Value = 168314  # was 160230
$ws.Range("G429").Value = 438  # was 585

$ws.Range("F430").Value = 164510  # was 152952
$ws.Range("G430").Value = 268  # was 249

$ws.Range("F431").Value = 162224  # was 135327
$ws.Range("G431").Value = 388  # was 288
